# "New Test Plans/Try TestCycle" — refresh the RMA references on the
# "RMA Details Maintenance Grid" sheet with a newly generated test cycle
# (RMA-87Y5 / RMA-PGV6), replacing the previous RMA-ZHZC cycle values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2
$ws.Range("E2").Value = "RMA-PGV6-001"
$ws.Range("F2").Value = "RMA-87Y5-1-1"
$ws.Range("J2").Value = "a7s5f000000xKZ5AAM"

# Row 3
$ws.Range("E3").Value = "RMA-PGV6-002"
$ws.Range("F3").Value = "RMA-87Y5-1-2"
$ws.Range("J3").Value = "a7s5f000000xKZ6AAM"

# Row 4
$ws.Range("E4").Value = "RMA-PGV6-003"
$ws.Range("F4").Value = "RMA-87Y5-1-3"
$ws.Range("J4").Value = "a7s5f000000xKZ7AAM"
